# Apply the "Inquiry" workbook edit:
#  - rename sheet from Inquiry_4 to Inquiry_1
#  - swap the two product names in rows 2 & 3 (biscuits <-> Sugar 5kg)
#  - update quantity / unit price / total amount figures for rows 2-4
#  - change row 4's status to "N/A" and add a new "Remarks" value in column G
#  - widen column F and give the new column G a width
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet/tab ---
$ws.Name = "Inquiry_1"

# --- Row 2: product swaps to "Sugar 5kg" with new quantities ---
$ws.Range("A2").Value = "Sugar 5kg"
$ws.Range("B2").Value = 50.0
$ws.Range("C2").Value = 200.0
$ws.Range("E2").Value = 10000.0

# --- Row 3: product swaps to "biscuits" with new quantities ---
$ws.Range("A3").Value = "biscuits"
$ws.Range("B3").Value = 30.0
$ws.Range("C3").Value = 150.0
$ws.Range("E3").Value = 4500.0

# --- Row 4: product becomes "Chicken", pricing + status updated ---
$ws.Range("A4").Value = "Chicken"
$ws.Range("C4").Value = 1000.0

# E4 and F4 change from the "Available" styling (s=3) to the plain (s=2)
# styling used by column D - copy formatting from D2 before setting values.
$ws.Range("D2").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = 0.0

$ws.Range("D2").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F4").Value = "N/A"

# New remarks column G4 - reuse the "Available" styling (s=3) that used
# to live on F4, by copying the format from F2 (still styled s=3).
$ws.Range("F2").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("G4").Value = "we only have 5 kg"

$excel.CutCopyMode = 0

# --- Column widths: widen F (Status) and size the new G (Remarks) column ---
$ws.Columns.Item(6).ColumnWidth = 18.6
$ws.Columns.Item(7).ColumnWidth = 21.6
